# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect refreshed counts captured at the newer data pull.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2214
    3 = 1677
    4 = 328
    5 = 1074
    6 = 711
    8 = 5779
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
